$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.608.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.141.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.67%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5265"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.92%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4555"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.85"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09139"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.184"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.136.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.875"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.147"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001171"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.010"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06704"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.008"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.346"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.720.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.389"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.371.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.54%  "

$ws.Range("E28").Value = "  +4.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.222"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1081"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.672"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.369"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.009"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.174"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.04%  "

$ws.Range("E38").Value = "  +2.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06929"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2336"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6984"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.271"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.350"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6460"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.59%  "

$ws.Range("E47").Value = "  +5.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.758"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.255"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("E51").Value = "  +2.37%  "
